$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new "2021" column (O) mirroring the existing "2020" column (N):
# first copy N4's formatting onto O4, then set its value.
$ws.Range("N4").Copy() | Out-Null
$ws.Range("O4").PasteSpecial(-4122) | Out-Null
$ws.Range("O4").Value = 2021

# Same for the data row: copy N5's formatting onto O5, then set its value.
$ws.Range("N5").Copy() | Out-Null
$ws.Range("O5").PasteSpecial(-4122) | Out-Null
$ws.Range("O5").Value = 1.5020015556876996

$excel.CutCopyMode = $false

# Update the saved selection to match the new active cell.
$ws.Range("Q5").Select() | Out-Null
